# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45205 to 45206 (bump the serial date by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
